$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data to reflect the latest scrape
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.600.09"
$ws.Range("E2").Value = "  -1.17%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.052.16"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.41"
$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.660"
$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.56"
$ws.Range("E8").Value = "  -7.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.67"
$ws.Range("E9").Value = "  +1.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.366"
$ws.Range("E10").Value = "  -3.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0752"
$ws.Range("E11").Value = "  -4.39%  "

$ws.Range("E12").Value = "  -3.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.964"
$ws.Range("E13").Value = "  +8.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.74"
$ws.Range("E14").Value = "  -4.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.353.57"
$ws.Range("E15").Value = "  -0.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.46"
$ws.Range("E16").Value = "  -4.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.069.25"
$ws.Range("E17").Value = "  +0.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.532.22"
$ws.Range("E18").Value = "  -1.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.12"
$ws.Range("E19").Value = "  -6.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.93"
$ws.Range("E20").Value = "  -2.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0858"
$ws.Range("E21").Value = "  -4.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.05"
$ws.Range("E22").Value = "  -0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.23"
$ws.Range("E23").Value = "  -4.01%  "

$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("E25").Value = "  -2.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.27"
$ws.Range("E26").Value = "  +4.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.41"
$ws.Range("E27").Value = "  -2.76%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.19"
$ws.Range("E28").Value = "  -10.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.06"
$ws.Range("E29").Value = "  -0.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.121"
$ws.Range("E30").Value = "  -2.31%  "

$ws.Range("E31").Value = "  +7.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.05"
$ws.Range("E32").Value = "  -8.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.47"
$ws.Range("E33").Value = "  -4.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0594"
$ws.Range("E34").Value = "  -3.99%  "

$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0866"
$ws.Range("E36").Value = "  +2.11%  "

$ws.Range("E37").Value = "  -0.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.23"
$ws.Range("E38").Value = "  -5.32%  "

$ws.Range("E39").Value = "  -6.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.03"
$ws.Range("E40").Value = "  -4.14%  "

$ws.Range("E41").Value = "  -5.40%  "

$ws.Range("E42").Value = "  -4.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.10"
$ws.Range("E43").Value = "  -5.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "94.35"
$ws.Range("E44").Value = "  -3.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0913"
$ws.Range("E45").Value = "  -5.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.414.64"
$ws.Range("E46").Value = "  +8.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.90"
$ws.Range("E47").Value = "  -6.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.47"
$ws.Range("E48").Value = "  +9.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.91"
$ws.Range("E49").Value = "  +1.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.26"
$ws.Range("E50").Value = "  -4.51%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.235.86"
$ws.Range("E51").Value = "  -0.45%  "
